$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 170
$ws1.Range("F5").Value = 4954
$ws1.Range("F9").Value = 542
$ws1.Range("F10").Value = 504
$ws1.Range("F11").Value = 29
$ws1.Range("F13").Value = 1371
$ws1.Range("F14").Value = 3561
$ws1.Range("F15").Value = 399
$ws1.Range("F16").Value = 131
$ws1.Range("F17").Value = 114
$ws1.Range("F18").Value = 77
$ws1.Range("F19").Value = 2605
$ws1.Range("F20").Value = 128
$ws1.Range("F22").Value = 36
$ws1.Range("F23").Value = 174
$ws1.Range("F24").Value = 43
$ws1.Range("F25").Value = 127
$ws1.Range("F26").Value = 55
$ws1.Range("F27").Value = 264
$ws1.Range("F28").Value = 45

# Sheet "全部类型" (sheet4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 170
$ws4.Range("F6").Value = 4954
$ws4.Range("F10").Value = 542
$ws4.Range("F11").Value = 504
$ws4.Range("F12").Value = 29
$ws4.Range("F14").Value = 1371
$ws4.Range("F15").Value = 3562
$ws4.Range("F16").Value = 399
$ws4.Range("F17").Value = 131
$ws4.Range("F18").Value = 114
$ws4.Range("F19").Value = 77
$ws4.Range("F20").Value = 2605
$ws4.Range("F21").Value = 128
$ws4.Range("F23").Value = 36
$ws4.Range("F24").Value = 174
$ws4.Range("F26").Value = 127
$ws4.Range("F27").Value = 55
$ws4.Range("F28").Value = 264
$ws4.Range("F29").Value = 45
